# Shiny Test Data.xlsx -- add deliberate "known errors" to the test data
# sheets so the QA/QC part has something to train against.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# plantData (sheet3): add percHerbPlant / numLeaves columns, plus two
# new rows that contain bad/odd values ("<5" text in a numeric column,
# a stray percentage) -- the "known errors" referenced by the commit.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("plantData")

$ws3.Range("C1").Value = "percHerbPlant"
$ws3.Range("D1").Value = "numLeaves"

$ws3.Range("C2").Value = 1
$ws3.Range("C3").Value = 34
$ws3.Range("C4").Value = 105

$ws3.Range("A5").Value = "plant"
$ws3.Range("B5").Value = 99
$ws3.Range("C5").Value = "<5"

$ws3.Range("A6").Value = "plant"
$ws3.Range("B6").Value = 99
$ws3.Range("C6").Value = 0.06
$ws3.Range("C6").NumberFormat = "0%"

# Column C sized to fit the "percHerbPlant" header (13 characters).
$ws3.Columns.Item(3).ColumnWidth = 12.166666666666666

$ws3.Range("D5").Select() | Out-Null

# ---------------------------------------------------------------------
# herbivoreData (sheet5): no data changes, just moves the live
# selection/active-tab elsewhere (see notes sheet below).
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("herbivoreData")
$ws5.Range("F8").Select() | Out-Null

# ---------------------------------------------------------------------
# notes (sheet7): drop the placeholder "notes / test dataset created"
# row entirely, and leave this sheet as the active one.
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("notes")
$ws7.Rows.Item(2).Delete() | Out-Null
$ws7.Range("B2").Select() | Out-Null
